$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 187.25
$ws.Range("I12").Value = 187.25
$ws.Range("K12").Value = 187.25
$ws.Range("M12").Value = -17.25
$ws.Range("H19").Value = 1513.25
$ws.Range("I19").Value = 1606.909
$ws.Range("J19").Value = 483
$ws.Range("K19").Value = 1606.909
$ws.Range("L19").Value = 483
$ws.Range("M19").Value = -1431.909
$ws.Range("N19").Value = -833
$ws.Range("H43").Value = 3932.6667
$ws.Range("I43").Value = 2400
$ws.Range("K43").Value = 2400
$ws.Range("M43").Value = -2331
$ws.Range("H47").Value = 12499
$ws.Range("I47").Value = 12499
$ws.Range("K47").Value = 12499
$ws.Range("M47").Value = -11527
$ws.Range("H86").Value = 1250
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 1250
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -18732
$ws.Range("H125").Value = 1454.75
$ws.Range("I125").Value = 1106.3334
$ws.Range("K125").Value = 9957.000599999999
$ws.Range("M125").Value = -7497.000599999999
$ws.Range("H138").Value = 2506.077
$ws.Range("I138").Value = 1298.5
$ws.Range("K138").Value = 3895.5
$ws.Range("M138").Value = 1244.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 700
$ws.Range("I132").Value = 700
$ws.Range("K132").Value = 2100
$ws.Range("M132").Value = 430

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 299
$ws.Range("I22").Value = 299
$ws.Range("K22").Value = 299
$ws.Range("M22").Value = -126
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H6").Value = 112183.336
$ws.Range("I6").Value = 126018.75
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 126018.75
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -125905.75
$ws.Range("N6").Value = -1726
$ws.Range("H7").Value = 760.069
$ws.Range("I7").Value = 80.5
$ws.Range("J7").Value = 937.34784
$ws.Range("K7").Value = 80.5
$ws.Range("L7").Value = 937.34784
$ws.Range("M7").Value = 32.5
$ws.Range("N7").Value = -1163.34784
$ws.Range("H22").Value = 212.625
$ws.Range("I22").Value = 240
$ws.Range("J22").Value = 167
$ws.Range("K22").Value = 240
$ws.Range("L22").Value = 167
$ws.Range("M22").Value = 110
$ws.Range("N22").Value = -867
$ws.Range("H23").Value = 16003333
$ws.Range("H25").Value = 526
$ws.Range("J25").Value = 526
$ws.Range("L25").Value = 526
$ws.Range("N25").Value = -874
$ws.Range("H27").Value = 16003333
$ws.Range("H31").Value = 4999.5
$ws.Range("I31").Value = 4999.5
$ws.Range("K31").Value = 4999.5
$ws.Range("M31").Value = -4704.5
$ws.Range("H34").Value = 4999.5
$ws.Range("I34").Value = 4999.5
$ws.Range("K34").Value = 4999.5
$ws.Range("M34").Value = -4797.5
$ws.Range("H43").Value = 16996.25
$ws.Range("J43").Value = 16996.25
$ws.Range("L43").Value = 16996.25
$ws.Range("N43").Value = -17364.25
$ws.Range("H101").Value = 16996.25
$ws.Range("J101").Value = 16996.25
$ws.Range("L101").Value = 16996.25
$ws.Range("N101").Value = -23486.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100106600
$ws.Range("I4").Value = 126998.125
$ws.Range("K4").Value = 380994.375
$ws.Range("M4").Value = -380882.375
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H41").Value = 1498
$ws.Range("I41").Value = 1498
$ws.Range("K41").Value = 4494
$ws.Range("M41").Value = -4156
$ws.Range("H55").Value = 2001.6666
$ws.Range("J55").Value = 2001.6666
$ws.Range("L55").Value = 6004.9998
$ws.Range("N55").Value = -6358.9998
$ws.Range("H68").Value = 493.33334
$ws.Range("I68").Value = 540
$ws.Range("K68").Value = 1620
$ws.Range("M68").Value = -809
$ws.Range("H71").Value = 493.33334
$ws.Range("I71").Value = 540
$ws.Range("K71").Value = 4860
$ws.Range("M71").Value = -804
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 500
$ws.Range("K97").Value = 1500
$ws.Range("M97").Value = -1004
$ws.Range("H103").Value = 2665.1428
$ws.Range("J103").Value = 2942.6667
$ws.Range("L103").Value = 8828.000100000001
$ws.Range("N103").Value = -10586.0001
$ws.Range("H121").Value = 1041.25
$ws.Range("I121").Value = 500
$ws.Range("K121").Value = 1500
$ws.Range("M121").Value = -190
$ws.Range("H122").Value = 994
$ws.Range("I122").Value = 994
$ws.Range("K122").Value = 8946
$ws.Range("M122").Value = -6496

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 35000
$ws.Range("I49").Value = 35000
$ws.Range("K49").Value = 35000
$ws.Range("M49").Value = -34816
$ws.Range("H132").Value = 4999.5
$ws.Range("I132").Value = 4999.5
$ws.Range("K132").Value = 14998.5
$ws.Range("M132").Value = -12468.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 675
$ws.Range("I16").Value = 675
$ws.Range("K16").Value = 675
$ws.Range("M16").Value = -505
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 5000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -4705
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 5000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -4893
$ws.Range("N27").ClearContents()
$ws.Range("H100").Value = 2421
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2421
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2421
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3503

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 14999
$ws.Range("I26").Value = 14999
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 14999
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -14706
$ws.Range("N26").ClearContents()
$ws.Range("H132").Value = 6454.2
$ws.Range("I132").Value = 6454.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19362.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -16832.6
$ws.Range("N132").ClearContents()

